$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D28").Value = "[논문리뷰] Control Barrier Functions : Theory and Applications"
$ws.Range("E28").Value = "https://ropiens.tistory.com/236"

$ws.Range("D44").Value = "[Book Review] 헤지펀드 열전"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/138"

$ws.Range("D51").Value = "[회사 용어] 일람 뜻"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%ED%9A%8C%EC%82%AC-%EC%9A%A9%EC%96%B4-%EC%9D%BC%EB%9E%8C-%EB%9C%BB"
